# "finished calculating and writing roto stats"
# Update the roto (fantasy) stats on Sheet2 with freshly calculated numbers,
# then leave the cursor/selection the way the author left it when they
# finished: Sheet1's selection collapsed to G1, and Sheet2 (the active
# sheet) selected at F5.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- updated roto stat values on Sheet2 ---
$ws2.Range("B2").Value = 83
$ws2.Range("B3").Value = 91
$ws2.Range("B4").Value = 83
$ws2.Range("F5").Value = 98
$ws2.Range("B11").Value = 92
$ws2.Range("B14").Value = 84

# --- selection / active sheet cleanup ---
$ws1.Range("G1").Select()
$ws2.Activate()
$ws2.Range("F5").Select()
